$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.15

$ws.Range("B3").Value = 1.54
$ws.Range("E3").Value = 1.33
$ws.Range("G3").Value = 0.65

$ws.Range("C4").Value = 1.45
$ws.Range("E4").Value = 1.23
$ws.Range("F4").Value = 1.11

$ws.Range("B5").Value = 1.58
$ws.Range("C5").Value = 1.34
$ws.Range("F5").Value = 1.06
$ws.Range("G5").Value = 0.77

$ws.Range("D6").Value = 1.5
$ws.Range("E6").Value = 1.32

$ws.Range("C7").Value = 2.17
$ws.Range("E7").Value = 1.88
